$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "27.111.68"
Set-TextCell "E2" "  -3.02%  "

Set-TextCell "D3" "1.869.00"
Set-TextCell "E3" "  -2.12%  "

Set-TextCell "D4" "1.001"
Set-TextCell "E4" "  +0.25%  "

Set-TextCell "D5" "307.38"
Set-TextCell "E5" "  -1.94%  "

Set-TextCell "E6" "  +0.25%  "

Set-TextCell "D7" "0.5059"
Set-TextCell "E7" "  +0.96%  "

Set-TextCell "D8" "0.3741"
Set-TextCell "E8" "  -2.17%  "

Set-TextCell "D9" "0.07153"
Set-TextCell "E9" "  -2.35%  "

Set-TextCell "D10" "0.8861"
Set-TextCell "E10" "  -3.00%  "

Set-TextCell "D11" "20.62"
Set-TextCell "E11" "  -2.73%  "

Set-TextCell "B12" "WrappedEther"
Set-TextCell "C12" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D12" "1.867.29"
Set-TextCell "E12" "  -0.68%  "

Set-TextCell "B13" "TRON"
Set-TextCell "C13" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell "D13" "0.07556"
Set-TextCell "E13" "  -1.61%  "

Set-TextCell "E14" "  -3.59%  "

Set-TextCell "D15" "89.14"
Set-TextCell "E15" "  -3.77%  "

Set-TextCell "D16" "1.002"
Set-TextCell "E16" "  +0.25%  "

Set-TextCell "E17" "  -3.40%  "

Set-TextCell "E18" "  -4.09%  "

Set-TextCell "D19" "1.000"
Set-TextCell "E19" "  +0.19%  "

Set-TextCell "D20" "27.162.26"
Set-TextCell "E20" "  -2.95%  "

Set-TextCell "D21" "5.070"
Set-TextCell "E21" "  -2.27%  "

Set-TextCell "D22" "2.078.96"
Set-TextCell "E22" "  -2.66%  "

Set-TextCell "E23" "  -2.88%  "

Set-TextCell "E24" "  -2.05%  "

Set-TextCell "D25" "150.62"
Set-TextCell "E25" "  -1.55%  "

Set-TextCell "D26" "1.838"
Set-TextCell "E26" "  -0.43%  "

Set-TextCell "D27" "17.98"
Set-TextCell "E27" "  -2.43%  "

Set-TextCell "D28" "2.090"
Set-TextCell "E28" "  -5.65%  "

Set-TextCell "D29" "112.66"
Set-TextCell "E29" "  -2.55%  "

Set-TextCell "D30" "4.746"
Set-TextCell "E30" "  -3.57%  "

Set-TextCell "D31" "4.690"
Set-TextCell "E31" "  -3.26%  "

Set-TextCell "D32" "0.09043"
Set-TextCell "E32" "  +0.16%  "

Set-TextCell "D33" "0.05126"
Set-TextCell "E33" "  -3.09%  "

Set-TextCell "D34" "3.091"
Set-TextCell "E34" "  -3.68%  "

Set-TextCell "D35" "0.7382"
Set-TextCell "E35" "  -5.32%  "

Set-TextCell "E36" "  -6.43%  "

Set-TextCell "D37" "0.02032"
Set-TextCell "E37" "  -2.68%  "

Set-TextCell "E38" "  -3.50%  "

Set-TextCell "D39" "3.045"
Set-TextCell "E39" "  -0.94%  "

Set-TextCell "D40" "1.080"
Set-TextCell "E40" "  -1.35%  "

Set-TextCell "D41" "0.5344"
Set-TextCell "E41" "  -4.08%  "

Set-TextCell "D42" "6.589"
Set-TextCell "E42" "  -4.47%  "

Set-TextCell "D43" "115.64"
Set-TextCell "E43" "  +2.13%  "

Set-TextCell "D44" "8.342"
Set-TextCell "E44" "  -2.30%  "

Set-TextCell "D45" "0.1474"
Set-TextCell "E45" "  -3.21%  "

Set-TextCell "E46" "  +0.25%  "

Set-TextCell "D47" "0.4629"
Set-TextCell "E47" "  -4.37%  "

Set-TextCell "E48" "  -6.05%  "

Set-TextCell "D49" "1.562"
Set-TextCell "E49" "  -4.90%  "

Set-TextCell "D50" "64.51"
Set-TextCell "E50" "  -4.69%  "

Set-TextCell "D51" "36.42"
Set-TextCell "E51" "  -2.03%  "
